$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "client/apllication/work"
$ws.Range("B13").Value = "work-template"
$ws.Range("C13").Value = "Work Template"

$ws.Range("A14").Value = "client/apllication/work"
$ws.Range("B14").Value = "safety-template"
$ws.Range("C14").Value = "Safety Template"

$ws.Range("C15").Select()
